# Adds a new data row (row 10) to Sheet1 that mirrors the existing
# "tytti@laulumuisto.fi" submission in row 4, but:
#   - uses a new "Vastausaika" of 29.1.2022
#   - leaves "Kompostoinnin vastuuhenkilon yhteystiedot" Etunimi/Sukunimi (D/E)
#     and the mirrored "Haltijan etunimi/sukunimi" (N/O) empty
#   - fills in "1. Kompostoria kayttavan rakennuksen tiedot" Etunimi/Sukunimi
#     (AW/AX) with Tytti / Tuntematon instead
#   - marks the row as "Virheellinen" (BT) instead of "Yksittainen"
#   - uses a later "Voimassaolopaiva" (BR)
# This mirrors the validator test-data change described in the commit
# message: a row whose vastuuhenkilo/kompostoija both lack etu- and
# sukunimi.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# NOTE: this runtime's PowerShell dialect does not bind named (-Foo bar)
# arguments correctly inside user functions, so positional parameters are
# used everywhere below.
function Set-CellFromTemplate($TemplateRow, $TemplateCol, $TargetRow, $TargetCol, $Value) {
    $src = $ws.Cells.Item($TemplateRow, $TemplateCol)
    $dst = $ws.Cells.Item($TargetRow, $TargetCol)
    $src.Copy()
    $dst.PasteSpecial($xlPasteFormats)
    $dst.Value = $Value
}

$srcRow = 4
$dstRow = 10

# column letter -> column number, value, (optional) template column override
$cells = @(
    @{ Col = 1;  Value = "29.1.2022" },                                   # A Vastausaika (new)
    @{ Col = 2;  Value = "Käsitelty kirjaamossa (tämän kohdan täyttää käsittelijä)" }, # B
    @{ Col = 3;  Value = "Muutos aikaisemman ilmoituksen tietoihin" },    # C
    # D (Etunimi) and E (Sukunimi) intentionally left blank
    @{ Col = 6;  Value = "0400123456" },                                  # F
    @{ Col = 8;  Value = "Hatustaheitettykatu 5B" },                      # H
    @{ Col = 9;  Value = "15100" },                                       # I
    @{ Col = 10; Value = "Lahti" },                                       # J
    @{ Col = 11; Value = "ASUNTO OY KAHDEN LAULUMUISTO" },                # K
    @{ Col = 13; Value = "Lahti" },                                       # M
    # N (Haltijan etunimi) and O (Haltijan sukunimi) intentionally left blank
    @{ Col = 16; Value = "HARJUKATU 44" },                                # P
    @{ Col = 17; Value = "15100" },                                       # Q
    @{ Col = 18; Value = "Lahti" },                                       # R
    @{ Col = 19; Value = "123456789A" },                                  # S
    @{ Col = 20; Value = "123456789A" },                                  # T
    @{ Col = 21; Value = "123456789A" },                                  # U
    @{ Col = 22; Value = "123456789A" },                                  # V
    @{ Col = 23; Value = "Kerrostalo" },                                  # W
    @{ Col = 25; Value = 5 },                                             # Y
    @{ Col = 26; Value = 1 },                                             # Z
    @{ Col = 27; Value = 160 },                                           # AA
    @{ Col = 30; Value = "SuperCompost" },                                # AD
    @{ Col = 32; Value = "Kyllä" },                                       # AF
    @{ Col = 33; Value = "Kyllä" },                                       # AG
    @{ Col = 34; Value = "Ei" },                                          # AH
    @{ Col = 35; Value = 12 },                                            # AI
    @{ Col = 36; Value = "Tammikuu" },                                    # AJ
    @{ Col = 37; Value = "Helmikuu" },                                    # AK
    @{ Col = 38; Value = "Maaliskuu" },                                   # AL
    @{ Col = 39; Value = "Huhtikuu" },                                    # AM
    @{ Col = 40; Value = "Toukokuu" },                                    # AN
    @{ Col = 41; Value = "Kesäkuu" },                                     # AO
    @{ Col = 42; Value = "Heinäkuu" },                                    # AP
    @{ Col = 43; Value = "Elokuu" },                                      # AQ
    @{ Col = 44; Value = "Syyskuu" },                                     # AR
    @{ Col = 45; Value = "Lokakuu" },                                     # AS
    @{ Col = 46; Value = "Marraskuu" },                                   # AT
    @{ Col = 47; Value = "Joulukuu" },                                    # AU
    @{ Col = 48; Value = "Kompostoria käyttää yksi rakennus, joka on ilmoitettu yllä Kompostorin sijainti -kohdassa" }, # AV
    @{ Col = 49; Value = "Tytti";       TemplateCol = 11 },               # AW (new name, style like a plain column)
    @{ Col = 50; Value = "Tuntematon" },                                  # AX (new name)
    @{ Col = 51; Value = "HARJUKATU 44" },                                # AY
    @{ Col = 52; Value = "15100" },                                       # AZ
    @{ Col = 53; Value = "Lahti" },                                       # BA
    @{ Col = 54; Value = "123456789A" },                                  # BB
    @{ Col = 55; Value = "123456789A" },                                  # BC
    @{ Col = 56; Value = "123456789A" },                                  # BD
    @{ Col = 57; Value = "Hyväksytty" },                                  # BE
    @{ Col = 58; Value = "Kerrostalo" },                                  # BF
    @{ Col = 69; Value = "Käsitelty" },                                   # BQ
    @{ Col = 70; Value = 46416 },                                         # BR
    @{ Col = 72; Value = "Virheellinen" }                                 # BT
)

foreach ($cell in $cells) {
    $templateCol = $cell.Col
    if ($cell.ContainsKey("TemplateCol")) {
        $templateCol = $cell.TemplateCol
    }
    Set-CellFromTemplate $srcRow $templateCol $dstRow $cell.Col $cell.Value
}

# G column: address/e-mail with a mailto hyperlink, same target as row 4.
$gCell = $ws.Cells.Item($dstRow, 7)
$gCell.Value = "tytti@laulumuisto.fi"
$ws.Hyperlinks.Add($gCell, "mailto:tytti@laulumuisto.fi", [Type]::Missing, [Type]::Missing, "tytti@laulumuisto.fi")
# Re-apply the same visual style used for the other e-mail cells (overrides
# the default hyperlink style Excel applies automatically).
$ws.Cells.Item($srcRow, 7).Copy()
$gCell.PasteSpecial($xlPasteFormats)

# Reflect the row being freshly added / selected, as happens when a user
# inserts & fills a new row at the bottom of the sheet.
$ws.Range("10:10").Select() | Out-Null

Write-Output "Row 10 added"
